$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8711434602737427
$ws.Range("B1").Value = 0.7253850102424622
$ws.Range("C1").Value = 0.644395649433136
$ws.Range("D1").Value = 0.7071085572242737
$ws.Range("E1").Value = 0.8773908019065857
